# Added Thieves Guild card images/list (Warriors/Mages/Shadow sections)
# to the "Units & Decks" sheet, below the existing Thalmor deck block.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New unique strings, written in original authoring order so the shared-string table
#     indices line up with the target workbook ---
$ws.Range("B246").Value = 'Deck: Thieves Guild'
$ws.Range("B248").Value = 'Maven Black Briar Leader card'
$ws.Range("B254").Value = 'red guard looter x2'
$ws.Range("B253").Value = 'corrupt guard skirmisher x2'
$ws.Range("B252").Value = 'corrupt guard x2'
$ws.Range("B283").Value = 'guild archer'
$ws.Range("B284").Value = 'thief'
$ws.Range("B285").Value = 'khajjit thief'
$ws.Range("B286").Value = 'fence'
$ws.Range("B287").Value = 'argonian thief'
$ws.Range("B288").Value = 'argonian archer'
$ws.Range("B267").Value = 'argonian picklocker'
$ws.Range("B255").Value = 'dirge'
$ws.Range("B256").Value = 'maul'
$ws.Range("B257").Value = 'hireling'
$ws.Range("B268").Value = 'Mauricio (hireling)'
$ws.Range("B289").Value = 'Karliah'
$ws.Range("B290").Value = 'Mercer Frey'
$ws.Range("H289").Value = 'Turn into nightingale'
$ws.Range("B258").Value = 'brynjolf'
$ws.Range("B293").Value = 'No Images:'
$ws.Range("B294").Value = 'Vipir The Fleet'
$ws.Range("B259").Value = 'gallus desidenius'
$ws.Range("D291").Value = 'spy'
$ws.Range("B291").Value = 'gulum Ei'
$ws.Range("B260").Value = 'gian the fist'
$ws.Range("B295").Value = 'vex'
$ws.Range("B296").Value = 'tonilia'
$ws.Range("B297").Value = 'delvin'
$ws.Range("B262").Value = 'no images'
$ws.Range("B263").Value = 'thrynn, rune'
$ws.Range("B298").Value = 'cynric endell'

# --- Remaining cells (numbers + reused existing strings) ---
$ws.Range("A244").Value = '##########################################################################################################'
$ws.Range("A250").Value = 'Warriors'
$ws.Range("A251").Value = 'No.'
$ws.Range("B251").Value = 'Name'
$ws.Range("C251").Value = 'Strength'
$ws.Range("D251").Value = 'Subtype'
$ws.Range("E251").Value = 'Created'
$ws.Range("F251").Value = 'Hero'
$ws.Range("G251").Value = 'Race'
$ws.Range("H251").Value = 'Ability'
$ws.Range("A252").Value = 1
$ws.Range("C252").Value = 1
$ws.Range("A253").Value = 2
$ws.Range("C253").Value = 2
$ws.Range("A254").Value = 3
$ws.Range("C254").Value = 3
$ws.Range("A255").Value = 4
$ws.Range("A256").Value = 5
$ws.Range("A257").Value = 6
$ws.Range("A258").Value = 7
$ws.Range("A259").Value = 8
$ws.Range("A260").Value = 9
$ws.Range("A265").Value = 'Mages'
$ws.Range("A266").Value = 'No.'
$ws.Range("B266").Value = 'Name'
$ws.Range("C266").Value = 'Strength'
$ws.Range("D266").Value = 'Subtype'
$ws.Range("E266").Value = 'Created'
$ws.Range("F266").Value = 'Hero'
$ws.Range("G266").Value = 'Race'
$ws.Range("H266").Value = 'Ability'
$ws.Range("A267").Value = 1
$ws.Range("A268").Value = 2
$ws.Range("C268").Value = 5
$ws.Range("D268").Value = 'spellsword'
$ws.Range("A269").Value = 3
$ws.Range("A270").Value = 4
$ws.Range("A271").Value = 5
$ws.Range("A272").Value = 6
$ws.Range("A273").Value = 7
$ws.Range("A274").Value = 8
$ws.Range("A281").Value = 'Shadow'
$ws.Range("A282").Value = 'No.'
$ws.Range("B282").Value = 'Name'
$ws.Range("C282").Value = 'Strength'
$ws.Range("D282").Value = 'Subtype'
$ws.Range("E282").Value = 'Created'
$ws.Range("F282").Value = 'Hero'
$ws.Range("G282").Value = 'Race'
$ws.Range("H282").Value = 'Ability'
$ws.Range("A283").Value = 1
$ws.Range("A284").Value = 2
$ws.Range("A285").Value = 3
$ws.Range("A286").Value = 4
$ws.Range("A287").Value = 5
$ws.Range("A288").Value = 6
$ws.Range("A289").Value = 7
$ws.Range("A290").Value = 8
$ws.Range("A291").Value = 9
$ws.Range("C291").Value = 2
$ws.Range("A294").Value = 9
$ws.Range("A295").Value = 10
$ws.Range("A296").Value = 11
$ws.Range("A297").Value = 12

# --- Restore the view state (scrolled/selected cell) the author left the sheet in ---
[void]$ws.Range("C258").Select()
$excel.ActiveWindow.ScrollRow = 243
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
